$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated run_time values (column C), rows 2-11
$cVals = @(0.5060012340545654,0.4659996032714844,0.5129795074462891,0.4409940242767334,0.4530165195465088,0.4250056743621826,0.5200014114379883,0.5470051765441895,0.3749980926513672,0.5160002708435059)

# Updated max_er + iter0..iter19 values (columns E:Y), rows 2-11
$eyVals = @(
    @(64.20996992708933,0.002802168546845917,0.002157148952229234,0.001719094163776776,0.001719094163776776,0.001687400868910901,0.00162145840520339,0.001565151783079396,0.001565151783079396,0.001565151783079396,0.001565151783079396,0.001565151783079396,0.001456809899573166,0.001387854942881566,0.001387854942881566,0.001387854942881566,0.001387854942881566,0.001301779425264658,0.001297807183691794,0.0012530583131119,0.00125165633386139),
    @(64.13189597493874,0.002698035534867245,0.002203862639851589,0.00200914947756402,0.00171310113755684,0.001695107453616733,0.001695107453616733,0.001501045378514589,0.001501045378514589,0.001501045378514589,0.001501045378514589,0.001501045378514589,0.001469831018379547,0.001403434284719813,0.001403434284719813,0.001329261552975374,0.001329261552975374,0.001272008431711285,0.001272008431711285,0.001262441782013301,0.001250134424462743),
    @(61.24760024835268,0.002802168546845917,0.002112673280831235,0.001927353959322397,0.001627564364483469,0.001627564364483469,0.001627564364483469,0.001627564364483469,0.001627564364483469,0.001528357097986461,0.001528357097986461,0.001460596993241226,0.001460596993241226,0.001460596993241226,0.001460596993241226,0.001414719050329521,0.001337199107350876,0.001264494293320052,0.001237442819587953,0.001236875086451654,0.001193910336225198),
    @(61.13960289489478,0.002650526418184921,0.001894014601547838,0.001894014601547838,0.001783419154449043,0.001643909720741786,0.001643909720741786,0.001471486339506744,0.001457956071455598,0.001457956071455598,0.001416958038063035,0.001416958038063035,0.001416958038063035,0.001414877401923546,0.001314883772445941,0.001314883772445941,0.001299811899019592,0.001273602954688063,0.001260809499500277,0.001222034487734817,0.001191805124656818),
    @(66.1323445942744,0.002802168546845917,0.002166852225781179,0.001909331412663204,0.001909331412663204,0.001733670217072961,0.001613736733779731,0.001613736733779731,0.001613254933269458,0.001479785924456891,0.001479785924456891,0.001405739925770969,0.001405739925770969,0.001405739925770969,0.001392315146607388,0.001392315146607388,0.00134690950750155,0.001346429404384962,0.00133061565597039,0.001302039871463197,0.001289129524254861),
    @(65.43603284732126,0.002678509318226792,0.002178762491041506,0.00198582904258473,0.001555915237822986,0.001555915237822986,0.001555915237822986,0.001450042690053008,0.001450042690053008,0.001450042690053008,0.001450042690053008,0.001442207344447255,0.001442207344447255,0.001399610632040571,0.00136781320630331,0.00136781320630331,0.001324930495314097,0.001320007277455397,0.001309923887229816,0.001302147993507378,0.001275556195854215),
    @(65.86575805266693,0.002802168546845917,0.002040791962604629,0.001916254492822424,0.001916254492822424,0.001771217275595455,0.001693358759261187,0.001629101600689632,0.001545920348642922,0.001489518175222668,0.001489518175222668,0.001489518175222668,0.001475399988666766,0.001460519087179714,0.001430583783935573,0.001333906787174775,0.001333906787174775,0.001333906787174775,0.001333906787174775,0.00131583305763215,0.001283932905510076),
    @(65.38283578494702,0.002725370669915938,0.002179206211326354,0.002179206211326354,0.001955369981386782,0.001924511131066474,0.001628018504854151,0.001628018504854151,0.001566393947538783,0.001566393947538783,0.001566393947538783,0.001566393947538783,0.001533139772623525,0.001506055260192137,0.001471990766812026,0.001419235968535583,0.001317910849410461,0.001317910849410461,0.001317910849410461,0.001286714615537128,0.001274519216080838),
    @(65.61485763809469,0.00265060364302398,0.002059624640169913,0.001999418176559129,0.00192062797526334,0.00188209440512784,0.001748883530750535,0.001748883530750535,0.001748883530750535,0.001607995401862045,0.001607995401862045,0.001529518667908104,0.001497084287544091,0.001492610363864472,0.001492610363864472,0.001492610363864472,0.00138319429027825,0.001355462137624548,0.001283106155628629,0.001283106155628629,0.001279042059222119),
    @(59.26464663557272,0.002636749212885127,0.002094206362705653,0.001873219184030335,0.0017397940553542,0.0017397940553542,0.001685824201066596,0.001643997036654718,0.00159137796740515,0.00159102006046732,0.001445786531311485,0.001445786531311485,0.001445786531311485,0.001309182996588277,0.001267543790593895,0.001256674819880473,0.001232935480015671,0.001214417968243985,0.001214417968243985,0.001176545416178536,0.001155256269699273)
)

for ($i = 0; $i -lt $cVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $rowVals = $eyVals[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 5
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
